# Prescription_sample.xlsx - anonymize STUDY_ID (col A) and MED_ORDER_ID (col C)
# for rows 2-79 on Sheet1, widen column A, and move the active selection to A13.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A (STUDY_ID): replace the real patient study IDs with small
# sequential placeholder numbers. The original values come in contiguous
# same-patient blocks, so set each block with one Range.Value assignment.
$ws.Range("A2:A2").Value = 1
$ws.Range("A3:A3").Value = 2
$ws.Range("A4:A9").Value = 3
$ws.Range("A10:A10").Value = 4
$ws.Range("A11:A42").Value = 5
$ws.Range("A43:A65").Value = 6
$ws.Range("A66:A79").Value = 7

# --- Column C (MED_ORDER_ID): replace the real order IDs with a simple
# sequential run starting at 11111 (row 2) through 11188 (row 79).
for ($r = 2; $r -le 79; $r++) {
    $ws.Cells.Item($r, 3).Value = 11109 + $r
}

# --- Column A width widened to fit the new header/content.
$ws.Columns("A").ColumnWidth = 25.83

# --- Move the selection/active cell to A13 (matches the saved view state).
$null = $ws.Range("A13").Select()
